$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.863.29'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '1.618.92'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.01'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.499'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('E8').Value = '  -1.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0615'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.32'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.26%  '
$ws.Range('E11').Value = '  -0.58%  '
$ws.Range('D12').Value = '1.843.14'
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('D13').Value = '1.615.99'
$ws.Range('E13').Value = '  -1.91%  '
$ws.Range('E14').Value = '  -2.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.522'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.05%  '
$ws.Range('D16').Value = '25.861.21'
$ws.Range('E16').Value = '  +0.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.33'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.60%  '
$ws.Range('D18').Value = '0.0₃0736'
$ws.Range('E18').Value = '  -2.61%  '
$ws.Range('E19').Value = '  +0.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '191.01'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.97%  '
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.46'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.01'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.134'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.21'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E27').Value = '  -3.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.67'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.89%  '
$ws.Range('E29').Value = '  -1.98%  '
$ws.Range('E31').Value = '  -2.47%  '
$ws.Range('E32').Value = '  -3.84%  '
$ws.Range('E33').Value = '  -4.77%  '
$ws.Range('E34').Value = '  -1.67%  '
$ws.Range('E35').Value = '  -2.83%  '
$ws.Range('D36').Value = '1.120.50'
$ws.Range('E36').Value = '  -0.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.832'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.38'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.94%  '
$ws.Range('E39').Value = '  -2.00%  '
$ws.Range('E40').Value = '  -4.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.11'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').Value = '1.753.73'
$ws.Range('E42').Value = '  -0.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.748'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.04'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.95%  '
$ws.Range('D45').Value = '0.0₆0113'
$ws.Range('E45').Value = '  +1.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '53.94'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.54%  '
$ws.Range('E47').Value = '  +0.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0521'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('E49').Value = '  -0.28%  '
$ws.Range('E50').Value = '  +0.40%  '
$ws.Range('E51').Value = '  -2.72%  '
